$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.947.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("E2").Style = "Normal"

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.805.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("E3").Style = "Normal"

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("E4").Style = "Normal"

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '703.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("E5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("E6").Style = "Normal"

# Row 7 - LidoStakedEther
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.804.61'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.51%  '
$ws.Range("E7").Style = "Normal"

# Row 8 - USDC
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E8").Style = "Normal"

# Row 9 - XRP
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("E9").Style = "Normal"

# Row 10 - Dogecoin
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("E10").Style = "Normal"

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.66'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.20%  '
$ws.Range("E11").Style = "Normal"

# Row 12 - Cardano
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.40%  '
$ws.Range("E12").Style = "Normal"

# Row 13 - ShibaInu
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("E13").Style = "Normal"

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("E14").Style = "Normal"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.444.88'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.61%  '
$ws.Range("E15").Style = "Normal"

# Row 16 - WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.805.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("E16").Style = "Normal"

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.905.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("E17").Style = "Normal"

# Row 18 - Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("E18").Style = "Normal"

# Row 19 - TRON
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("E19").Style = "Normal"

# Row 20 - Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("E20").Style = "Normal"

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '499.24'
$ws.Range("D21").Style = "Normal"

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("E22").Style = "Normal"

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("E23").Style = "Normal"

# Row 24 - Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.09%  '
$ws.Range("E24").Style = "Normal"

# Row 25 - PEPE
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.92%  '
$ws.Range("E25").Style = "Normal"

# Row 26 - WrappedeETH
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.955.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("E26").Style = "Normal"

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("E27").Style = "Normal"

# Row 28 - RenderToken
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.96%  '
$ws.Range("E28").Style = "Normal"

# Row 29 - Dai
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E29").Style = "Normal"

# Row 30 - Fetch.AI
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.09%  '
$ws.Range("E30").Style = "Normal"

# Row 31 - PancakeSwap
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.28%  '
$ws.Range("E31").Style = "Normal"

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("E32").Style = "Normal"

# Row 33 - NEARProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.49%  '
$ws.Range("E33").Style = "Normal"

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.99'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.20%  '
$ws.Range("E34").Style = "Normal"

# Row 35 - Kaspa
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.30%  '
$ws.Range("E35").Style = "Normal"

# Row 36 - Binance-PegBSC-USD -> RenzoRestakedETH
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.770.70'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E36").Style = "Normal"

# Row 37 - RenzoRestakedETH -> Binance-PegBSC-USD
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("E37").Style = "Normal"

# Row 38 - Aptos
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("E38").Style = "Normal"

# Row 39 - Hedera
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("E39").Style = "Normal"

# Row 40 - Stacks
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("E40").Style = "Normal"

# Row 41 - Filecoin
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E41").Style = "Normal"

# Row 42 - Mantle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.25%  '
$ws.Range("E42").Style = "Normal"

# Row 43 - dogwifhat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.24'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.25%  '
$ws.Range("E43").Style = "Normal"

# Row 45 - FirstDigitalUSD
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E45").Style = "Normal"

# Row 46 - Monero
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '166.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("E46").Style = "Normal"

# Row 47 - FLOKI
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000315'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.44%  '
$ws.Range("E47").Style = "Normal"

# Row 48 - OKB
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("E48").Style = "Normal"

# Row 49 - Bittensor
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("E49").Style = "Normal"

# Row 50 - Cosmos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.60%  '
$ws.Range("E50").Style = "Normal"

# Row 51 - ONDO
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.81%  '
$ws.Range("E51").Style = "Normal"
